$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$questionsText = @"
questions = [
    {
        "title": "How would you make multiple large files suitable for online transfer, without reducing the quality?",
        "ques_type": 2,
        "options": [
            "Open the files and remove content that requires large amounts of storage space (e.g. high quality pictures).",
            "Put them in a compressed (zip) folder.",
            "Upload the files to an online file compressor.",
            "Transfer all the files separately."
        ],
        "score": "Put them in a compressed (zip) folder."
    },
    {
        "title": "What shortcut keys would you use to cut a paragraph from a Word document and paste it in an e-mail?",
        "ques_type": 2,
        "options": [
            "CTRL + C and CTRL + V",
            "CTRL + S and CTRL + F",
            "CTRL + X and CTRL + V",
            "CTRL + Q and CTRL + X"
        ],
        "score": "CTRL + X and CTRL + V"
    },
    {
        "title": "How can you define a cursor?",
        "ques_type": 2,
        "options": [
            "A graphic pointer on the screen that shows the user where the mouse is.",
            "A virus on a computer that removes all internal memory.",
            "A symbol connected to a file that indicates the type of file (e.g. Word-document, PDF).",
            "A hacker that attempts to steal personal files from a computer."
        ],
        "score": "A graphic pointer on the screen that shows the user where the mouse is."
    },
    {
        "title": "You created an invoice in Microsoft Word. You would like to send the invoice to a customer. However, the customer is not allowed to make any changes to the invoice. How would you make sure the customer cannot make any changes?",
        "ques_type": 2,
        "options": [
            "Save the invoice as a PDF file.",
            "Save the invoice as a Fixed Text Document.",
            "Save the invoice as an RTF file.",
            "Save the invoice as a Word Template."
        ],
        "score": "Save the invoice as a PDF file."
    }
]
"@

# Clear the old A2 cell (it held the old shared-string question blob)
$ws.Range("A2").ClearContents()

# Remove the bold/bordered/centered style that used to be on A1
$ws.Range("A1").ClearFormats()

# Write the new pretty-printed JSON text into A1
$ws.Range("A1").Value = $questionsText

# Let Excel re-measure row 1's height for the new content, then drop back
# to the sheet's default (un-set any custom height) so the row matches
# the workbook's normal formatting.
$ws.Rows(1).AutoFit()
